# Refresh the crypto "Price" (column D) and "Volume(1h)" (column E) readings
# for the symbol list, matching the latest scrape (5-1-2023, hour 22).
# Values must stay as text (they include trailing "%", thousands separators,
# and fixed decimal padding like "256.90"), so each cell is forced to text
# format before the write and then has that formatting stripped again so no
# lingering style/number-format is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "256.90"
    "E2" = "-0.26%"
    "E3" = "-0.79%"
    "D4" = "4.619"
    "E4" = "-11.40%"
    "D5" = "0.05881"
    "E5" = "-0.66%"
    "D6" = "6.626"
    "E6" = "-0.57%"
    "D7" = "0.8579"
    "E7" = "-0.54%"
    "D8" = "0.9409"
    "E8" = "-8.45%"
    "D9" = "0.1409"
    "E9" = "-0.46%"
    "D10" = "0.04161"
    "E10" = "17.14%"
    "D11" = "0.07102"
    "E11" = "-1.09%"
    "D12" = "0.03183"
    "E12" = "1.68%"
    "D13" = "0.09153"
    "E13" = "-0.86%"
    "D14" = "0.001547"
    "E14" = "-0.94%"
    "D15" = "0.0006040"
    "E15" = "-0.74%"
    "D16" = "0.006219"
    "E16" = "8.86%"
    "D17" = "3.528"
    "E17" = "0.92%"
    "D18" = "3.206"
    "E18" = "-1.89%"
    "E19" = "0.64%"
    "D20" = "0.3052"
    "E20" = "-2.89%"
    "E21" = "-0.42%"
    "D22" = "3.833"
    "E22" = "8.91%"
    "D23" = "0.04239"
    "E23" = "1.33%"
    "D24" = "0.001223"
    "E24" = "0.44%"
    "D25" = "0.004294"
    "E25" = "-4.93%"
    "E26" = "0.04%"
    "E27" = "30.58%"
    "D40" = "0.03831"
    "E40" = "0.45%"
    "D41" = "0.006230"
    "E41" = "55.68%"
    "D42" = "0.1102"
    "E42" = "-0.32%"
    "E43" = "-6.74%"
    "D44" = "0.01147"
    "E44" = "6.98%"
    "D45" = "0.00005472"
    "E45" = "0.85%"
    "E46" = "0.08%"
    "D47" = "0.05121"
    "E47" = "-53.02%"
    "D48" = "0.1968"
    "E48" = "8,709.50%"
    "E49" = "0.08%"
    "E50" = "0.08%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
